$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.544047297731126
$ws.Range("C2").Value = 0.5520118787769169
$ws.Range("D2").Value = 0.4161555526239861
$ws.Range("E2").Value = 0.6451011956460677
$ws.Range("F2").Value = 0.3588196280695601
$ws.Range("G2").Value = 15

$ws.Range("B3").Value = 0.3323214535345447
$ws.Range("C3").Value = 0.3471491748320992
$ws.Range("D3").Value = 0.171833097400013
$ws.Range("E3").Value = 0.4145275592768387
$ws.Range("F3").Value = 0.2571347573803706
$ws.Range("G3").Value = 14

$ws.Range("B4").Value = 0.2658657847736177
$ws.Range("C4").Value = 0.3120893012315514
$ws.Range("D4").Value = 0.1488418498313651
$ws.Range("E4").Value = 0.3858002719430937
$ws.Range("F4").Value = 0.2909816784230805
$ws.Range("G4").Value = 13

$ws.Range("B5").Value = 0.3912245664449679
$ws.Range("C5").Value = 0.4218923612313186
$ws.Range("D5").Value = 0.2401193300839004
$ws.Range("E5").Value = 0.4900197241784257
$ws.Range("F5").Value = 0.3081841280094164
$ws.Range("G5").Value = 12

$ws.Range("B6").Value = 0.4189114456561581
$ws.Range("C6").Value = 0.4471377096428504
$ws.Range("D6").Value = 0.2554094045634564
$ws.Range("E6").Value = 0.5053804552646021
$ws.Range("F6").Value = 0.2965044110766253
$ws.Range("G6").Value = 11

$ws.Range("B7").Value = 0.3623555131786818
$ws.Range("C7").Value = 0.3958820287295867
$ws.Range("D7").Value = 0.197457793314669
$ws.Range("E7").Value = 0.4443622320974961
$ws.Range("F7").Value = 0.271121693430343
$ws.Range("G7").Value = 10

$ws.Range("B8").Value = 0.356182384604201
$ws.Range("C8").Value = 0.3864723693532894
$ws.Range("D8").Value = 0.1809534852127003
$ws.Range("E8").Value = 0.4253862776497384
$ws.Range("F8").Value = 0.2466749751680558
$ws.Range("G8").Value = 9

$ws.Range("B9").Value = 0.3694669697261536
$ws.Range("C9").Value = 0.3995817503120815
$ws.Range("D9").Value = 0.1961884184341418
$ws.Range("E9").Value = 0.4429316182371064
$ws.Range("F9").Value = 0.2611678753281997
$ws.Range("G9").Value = 8

$ws.Range("B10").Value = 0.3583150241653129
$ws.Range("C10").Value = 0.3939248673888729
$ws.Range("D10").Value = 0.1944778797009319
$ws.Range("E10").Value = 0.4409964622317643
$ws.Range("F10").Value = 0.277674138187313
$ws.Range("G10").Value = 7

$ws.Range("B11").Value = 0.3250039527724223
$ws.Range("C11").Value = 0.3520830062954648
$ws.Range("D11").Value = 0.1591455960645969
$ws.Range("E11").Value = 0.3989305654679732
$ws.Range("F11").Value = 0.25341987312813
$ws.Range("G11").Value = 6

